$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("I3").Value = 6
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 13

# Row 4
$ws.Range("O4").Value = 1.53
$ws.Range("P4").Value = 2.5
$ws.Range("Q4").Value = 2.6
$ws.Range("R4").Value = 1.48

# Row 5
$ws.Range("G5").Value = 1.5
$ws.Range("H5").Value = 3.8
$ws.Range("I5").Value = 7.5
$ws.Range("J5").Value = 2.1

# Row 6
$ws.Range("G6").Value = 2.2
$ws.Range("I6").Value = 3.2
$ws.Range("J6").Value = 2.88
$ws.Range("N6").Value = 9.5
$ws.Range("X6").Value = 10
$ws.Range("AJ6").Value = 12
$ws.Range("AO6").Value = 12
$ws.Range("AY6").Value = 19

# Row 16
$ws.Range("G16").Value = 3.35
$ws.Range("H16").Value = 4.45
$ws.Range("I16").Value = 1.7
$ws.Range("J16").Value = 3.45
$ws.Range("K16").Value = 2.82
$ws.Range("L16").Value = 2.07
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("O16").Value = 1.07
$ws.Range("P16").Value = 7.5
$ws.Range("Q16").Value = 1.16
$ws.Range("R16").Value = 3.94
$ws.Range("S16").Value = 1.12
$ws.Range("T16").Value = 5.3
$ws.Range("U16").Value = 1.24
$ws.Range("V16").Value = 3.87
$ws.Range("W16").Value = 24
$ws.Range("X16").Value = 27
$ws.Range("Y16").Value = 12.5
$ws.Range("Z16").Value = 45
$ws.Range("AB16").Value = 16.5
$ws.Range("AC16").Value = 37
$ws.Range("AD16").Value = 10.75
$ws.Range("AE16").Value = 10.25
$ws.Range("AF16").Value = 19.5
$ws.Range("AG16").Value = 60
$ws.Range("AH16").Value = 17
$ws.Range("AI16").Value = 14
$ws.Range("AJ16").Value = 8.75
$ws.Range("AK16").Value = 16.5
$ws.Range("AM16").Value = 12
$ws.Range("AN16").Value = 7
$ws.Range("AP16").Value = 14.5
$ws.Range("AQ16").Value = 55
$ws.Range("AR16").Value = 50
$ws.Range("AS16").Value = 80
$ws.Range("AT16").Value = 5.6
$ws.Range("AV16").Value = 24
$ws.Range("AW16").Value = 200
$ws.Range("AX16").Value = 4.9
$ws.Range("AZ16").Value = 10.25
$ws.Range("BA16").Value = 19.5
$ws.Range("BB16").Value = 25
$ws.Range("BC16").Value = 60

# Row 19
$ws.Range("G19").Value = 1.29
$ws.Range("H19").Value = 5.2
$ws.Range("I19").Value = 9.25
$ws.Range("J19").Value = 1.65
$ws.Range("K19").Value = 2.65
$ws.Range("L19").Value = 7.6
$ws.Range("P19").Value = 4.45
$ws.Range("Q19").Value = 1.52
$ws.Range("R19").Value = 2.37
$ws.Range("T19").Value = 3.5
$ws.Range("U19").Value = 1.9
$ws.Range("W19").Value = 8.25
$ws.Range("Y19").Value = 8.75
$ws.Range("AA19").Value = 10.25
$ws.Range("AD19").Value = 10.5
$ws.Range("AG19").Value = 700
$ws.Range("AJ19").Value = 28
$ws.Range("AL19").Value = 110
$ws.Range("AM19").Value = 80
$ws.Range("AO19").Value = 5.5
$ws.Range("AQ19").Value = 13
$ws.Range("AR19").Value = 32
$ws.Range("AS19").Value = 150
$ws.Range("AT19").Value = 3.5
$ws.Range("AU19").Value = 8.5
$ws.Range("AV19").Value = 70
$ws.Range("AX19").Value = 10
$ws.Range("AZ19").Value = 45
$ws.Range("BA19").Value = 400
$ws.Range("BB19").Value = 350
